$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 152.22223
$ws.Range("I6").Value = 152.22223
$ws.Range("K6").Value = 456.66669
$ws.Range("M6").Value = -344.66669
$ws.Range("H9").Value = 388.5
$ws.Range("I9").Value = 334.33334
$ws.Range("K9").Value = 334.33334
$ws.Range("M9").Value = -165.33334
$ws.Range("H12").Value = 921.5714
$ws.Range("I12").Value = 824.5
$ws.Range("J12").Value = 960.4
$ws.Range("K12").Value = 824.5
$ws.Range("L12").Value = 960.4
$ws.Range("M12").Value = -654.5
$ws.Range("N12").Value = -1300.4
$ws.Range("H21").Value = 21666.666
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 21666.666
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1500
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -2062
$ws.Range("H38").Value = 177.5
$ws.Range("I38").Value = 33.2
$ws.Range("J38").Value = 899
$ws.Range("K38").Value = 99.60000000000001
$ws.Range("L38").Value = 2697
$ws.Range("M38").Value = 272.4
$ws.Range("N38").Value = -3441
$ws.Range("H80").Value = 1879.8
$ws.Range("I80").Value = 899
$ws.Range("J80").Value = 2125
$ws.Range("K80").Value = 2697
$ws.Range("L80").Value = 6375
$ws.Range("M80").Value = -1699
$ws.Range("N80").Value = -8371
$ws.Range("H83").Value = 1879.8
$ws.Range("I83").Value = 899
$ws.Range("J83").Value = 2125
$ws.Range("K83").Value = 8091
$ws.Range("L83").Value = 19125
$ws.Range("M83").Value = -3099
$ws.Range("N83").Value = -29109
$ws.Range("H112").Value = 1827.2727
$ws.Range("H132").Value = 2279.4614
$ws.Range("I132").Value = 2499.5454
$ws.Range("K132").Value = 7498.6362
$ws.Range("M132").Value = -4968.6362
$ws.Range("H138").Value = 6429.64
$ws.Range("J138").Value = 6760.6313
$ws.Range("L138").Value = 20281.8939
$ws.Range("N138").Value = -30561.8939

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 2067.75
$ws.Range("I28").Value = 2067.75
$ws.Range("K28").Value = 2067.75
$ws.Range("M28").Value = -1875.75
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H32").Value = 8981.532999999999
$ws.Range("I32").Value = 8981.532999999999
$ws.Range("K32").Value = 8981.532999999999
$ws.Range("M32").Value = -8694.532999999999
$ws.Range("H99").Value = 2067.75
$ws.Range("I99").Value = 2067.75
$ws.Range("K99").Value = 2067.75
$ws.Range("M99").Value = 927.25
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 2492.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 38000
$ws.Range("J35").Value = 38000
$ws.Range("L35").Value = 38000
$ws.Range("N35").Value = -38620
$ws.Range("H47").Value = 199999
$ws.Range("J47").Value = 199999
$ws.Range("L47").Value = 199999
$ws.Range("N47").Value = -201039
$ws.Range("H97").Value = 13622.25
$ws.Range("I97").Value = 13622.25
$ws.Range("K97").Value = 13622.25
$ws.Range("M97").Value = -12631.25
$ws.Range("H102").Value = 24923.6
$ws.Range("I102").Value = 18751.5
$ws.Range("K102").Value = 18751.5
$ws.Range("M102").Value = -15506.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 929.5714
$ws.Range("I4").Value = 1036.6666
$ws.Range("J4").Value = 849.25
$ws.Range("K4").Value = 1036.6666
$ws.Range("L4").Value = 849.25
$ws.Range("M4").Value = -924.6666
$ws.Range("N4").Value = -1073.25
$ws.Range("H31").Value = 8571.286
$ws.Range("J31").Value = 9333.166999999999
$ws.Range("L31").Value = 9333.166999999999
$ws.Range("N31").Value = -9923.166999999999
$ws.Range("H34").Value = 8571.286
$ws.Range("J34").Value = 9333.166999999999
$ws.Range("L34").Value = 9333.166999999999
$ws.Range("N34").Value = -9737.166999999999
$ws.Range("H132").Value = 1024.75
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 3000
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1034.625
$ws.Range("I98").Value = 399
$ws.Range("J98").Value = 1416
$ws.Range("K98").Value = 1197
$ws.Range("L98").Value = 4248
$ws.Range("M98").Value = 301
$ws.Range("N98").Value = -7244
$ws.Range("H113").Value = 983.3333
$ws.Range("J113").Value = 975
$ws.Range("L113").Value = 2925
$ws.Range("N113").Value = -7265

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 33004
$ws.Range("I99").Value = 19246.5
$ws.Range("K99").Value = 19246.5
$ws.Range("M99").Value = -17000.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H98").Value = 25000
$ws.Range("J98").Value = 25000
$ws.Range("L98").Value = 25000
$ws.Range("N98").Value = -30990
$ws.Range("H122").Value = 5375.5
$ws.Range("I122").Value = 4606.25
$ws.Range("K122").Value = 13818.75
$ws.Range("M122").Value = -11368.75
$ws.Range("H132").Value = 3594.8948
$ws.Range("I132").Value = 3860.6
$ws.Range("K132").Value = 11581.8
$ws.Range("M132").Value = -9051.799999999999
$ws.Range("H136").Value = 1601.5
$ws.Range("I136").Value = 1601.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4804.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2254.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H23").Value = 10012470
$ws.Range("I23").Value = 599
$ws.Range("J23").Value = 12515438
$ws.Range("K23").Value = 599
$ws.Range("L23").Value = 12515438
$ws.Range("M23").Value = -370
$ws.Range("N23").Value = -12515896
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H100").Value = 3720.2
$ws.Range("I100").Value = 3251.2856
$ws.Range("K100").Value = 6502.5712
$ws.Range("M100").Value = -5961.5712
$ws.Range("H106").Value = 50000
$ws.Range("I106").Value = 50000
$ws.Range("K106").Value = 50000
$ws.Range("M106").Value = -48738
$ws.Range("H109").Value = 30000
$ws.Range("I109").Value = 30000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 30000
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -28613
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 400
$ws.Range("I132").Value = 400
$ws.Range("K132").Value = 1200
$ws.Range("M132").Value = 1330
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -110119
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
